# maj perk data use
# - Normalize the "Rarity" column (C) values to lowercase (Epic/Rare/Uncommon -> epic/rare/uncommon)
# - Remove the "READ FIRST" instructions textbox/shape from the sheet
# - Move the active selection to F12 (matches the saved cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lower-case every "Rarity" value found in column C (data rows start at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = 13 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value = $val.ToString().ToLower()
    }
}

# Delete any shapes on the sheet (removes the "READ FIRST" textbox / drawing part)
$shapeCount = $ws.Shapes.Count
for ($i = $shapeCount; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# Restore the saved selection location
$ws.Range("F12").Select()
